$d = $word.ActiveDocument

# --- 1. First paragraph: append "  " to the existing sentence, then three
#        red-colored runs forming "(This is a change – Version for main branch)"
$p1 = $d.Paragraphs.Item(1)
$r  = $p1.Range
# Range covering just the existing sentence text, excluding the paragraph mark.
$body = $d.Range($r.Start, $r.End - 1)
$body.Collapse(0)
$body.InsertAfter("  ")

$body.Collapse(0)
$body.InsertAfter("(This is a change – Ve")
$body.Font.Color = 255

$body.Collapse(0)
$body.InsertAfter("rsion for main branch")
$body.Font.Color = 255

$body.Collapse(0)
$body.InsertAfter(")")
$body.Font.Color = 255

# --- 2. Third paragraph (the empty Menlo-styled one) becomes a bare, empty
#        paragraph with no paragraph/run properties at all.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Delete()

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.Text = "`r"

# --- 3. Drop the handful of built-in/custom styles that are not referenced
#        anywhere in the document body (unused-style cleanup). Deleted in
#        reverse definition order so earlier removals never shift the index
#        of a not-yet-processed style.
$unusedStyles = @(
    "FollowedHyperlink",
    "c-txt",
    "Heading1Char",
    "podcast-toolssubscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading4Char",
    "Heading2Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading4",
    "Heading2",
    "Heading1"
)
foreach ($styleName in $unusedStyles) {
    $style = $d.Styles.Item($styleName)
    $style.Delete()
}
